# New model Sept 21
#
# Insert two new tag rows ("dom" and "dom-interview", with no recorded
# points yet) into the sorted N:P tag-summary table on Sheet1, right
# after "datetime" and before "easy" (their correct alphabetical
# slot). Every row of the summary table from "easy" onward therefore
# shifts down by two rows. The A:I table to the left is left
# completely untouched - only columns N, O and P are affected.
#
# Rows 2-7 of the table ("arrays" .. "datetime") are unaffected by the
# insert and are left alone; only rows 8-28 need to be rewritten.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final (tag, points) content for rows 8..28 of the N:P summary table
# after the insert.
$changedRows = @(
    @{ Row = 8;  Tag = "dom";                 Points = $null },
    @{ Row = 9;  Tag = "dom-interview";       Points = $null },
    @{ Row = 10; Tag = "easy";                Points = 5 },
    @{ Row = 11; Tag = "easy-medium";         Points = 10 },
    @{ Row = 12; Tag = "es6";                 Points = 4 },
    @{ Row = 13; Tag = "functions / methods"; Points = 4 },
    @{ Row = 14; Tag = "hard";                Points = 20 },
    @{ Row = 15; Tag = "hashes";              Points = 3 },
    @{ Row = 16; Tag = "hashes / js objects"; Points = 4 },
    @{ Row = 17; Tag = "HR";                  Points = 1 },
    @{ Row = 18; Tag = "javascript";          Points = 1 },
    @{ Row = 19; Tag = "loops";               Points = 1 },
    @{ Row = 20; Tag = "medium";              Points = 15 },
    @{ Row = 21; Tag = "methods";             Points = 4 },
    @{ Row = 22; Tag = "numbers";             Points = 3 },
    @{ Row = 23; Tag = "oops";                Points = 3 },
    @{ Row = 24; Tag = "prep";                Points = 1 },
    @{ Row = 25; Tag = "ruby";                Points = 1 },
    @{ Row = 26; Tag = "strings";             Points = 2 },
    @{ Row = 27; Tag = "variables";           Points = 3 },
    @{ Row = 28; Tag = "w3r";                 Points = 1 }
)

# Rows that must carry the alternating gray-fill banding (style index
# 4 in the original workbook). Before the edit this landed on rows 8
# and 11 ("easy" and "functions / methods"); after inserting the two
# blank rows, that same banding now belongs to rows 10 and 13.
$styledRows = @(10, 13)
$unstyledRows = @(8, 9, 11)

# Move the banding so it lines up with the new row positions: copy the
# format from row 2 (which still keeps its banding) onto the rows that
# now need it, and reset the rows that no longer need it back to the
# default "Normal" style.
$ws.Range("N2:P2").Copy() | Out-Null
foreach ($r in $styledRows) {
    $ws.Range("N" + $r + ":P" + $r).PasteSpecial(-4122) | Out-Null
}
$excel.CutCopyMode = 0

foreach ($r in $unstyledRows) {
    $ws.Range("N" + $r + ":P" + $r).Style = "Normal"
}

# Write out the final tag / points / formula content for rows 8-28.
foreach ($item in $changedRows) {
    $r = $item.Row
    $ws.Range("N" + $r).Value = $item.Tag
    if ($null -eq $item.Points) {
        $ws.Range("O" + $r).Value = ""
        $ws.Range("P" + $r).Formula = ""
    } else {
        $ws.Range("O" + $r).Value = $item.Points
        $ws.Range("P" + $r).Formula = "=O" + $r + "*5"
    }
}

# Update the cached sort-state range to reflect the new extent of the
# sorted table (it used to cover N2:O26, now it covers N2:O28).
$sortObj = $ws.Sort
$sortObj.SortFields.Clear() | Out-Null
$sortObj.SortFields.Add($ws.Range("N1")) | Out-Null
$sortObj.SetRange($ws.Range("N2:O28")) | Out-Null
$sortObj.Header = -4142
$sortObj.Apply() | Out-Null

# The user's selection ended up on O8 after entering the new data.
$ws.Range("O8").Select() | Out-Null
